$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-02-04T10:58:36+00:00"

# --- Elements sheet: add a new mapping row for the radiopharmaceutical
#     administration entry, right after the existing "quantiteExposition"
#     entry (row 9). ---
$elements = $wb.Worksheets.Item("Elements")

# Clone row 9 (values + style) into the new row 10, then overwrite just the
# cells that differ for the new entry.
$elements.Range("A9:AJ9").Copy($elements.Range("A10:AJ10"))

$elements.Range("A10").Value = "fr-lm-exposition-radiations.entree.administrationRadiopharmaceutique"
$elements.Range("B10").Value = "fr-lm-exposition-radiations.entree.administrationRadiopharmaceutique"
$elements.Range("G10").Value = "1"
$elements.Range("K10").Value = "https://interop.esante.gouv.fr/ig/document/core/StructureDefinition/fr-lm-administration-produit-de-sante`n"
$elements.Range("L10").Value = "Entrée administration des produits radiopharmaceutiques"
$elements.Range("M10").Value = "Entrée administration des produits radiopharmaceutiques"
$elements.Range("AF10").Value = "fr-lm-exposition-radiations.entree.administrationRadiopharmaceutique"
$elements.Range("AH10").Value = "1"
